$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the pt_max column (F2:F12) values from 60 to 59
$ws.Range("F2:F12").Value = 59
